$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date value (45204 = 2023-10-05) that needs
# to be bumped by one day (45205 = 2023-10-06) for every data row (2-78).
$ws.Range("C2:C78").Value = 45205
